$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet plumbing: turn the current "总计" sheet into the new "2022-Q1"
#    sheet (it keeps sheetId=2) and create a brand-new "总计" sheet (a copy
#    of the old one, picking up sheetId=3) so the final tab order is
#    2021-Q4, 2022-Q1, 总计 - matching the target workbook.xml.
# ---------------------------------------------------------------------------
$src = $wb.Worksheets.Item("总计")
$src.Copy($null, $src)

$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Move($null, $wb.Worksheets.Item("2021-Q4"))

$total = $wb.Worksheets.Item("总计 (2)")
$total.Name = "总计"

# ---------------------------------------------------------------------------
# 2) "2022-Q1" sheet: replace the old "总计" summary content with the
#    per-fund holder detail table (columns A-H).
# ---------------------------------------------------------------------------

# Propagate the existing header style (currently only on B1:D1) across the
# new columns E1:H1 before filling the text in.
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2 (index 0) keeps A2's existing style; just overwrite values. The
# numeric-looking text columns (D-G) are quote-prefixed so they stay text
# ("2.31") instead of being auto-coerced to a number (2.31).
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'090019"
$q1.Range("C2").Value = "大成景恒混合A"
$q1.Range("D2").Value = "'2.31"
$q1.Range("E2").Value = "'93.51"
$q1.Range("F2").Value = "'2.08"
$q1.Range("G2").Value = "'0.0480"
$q1.Range("H2").Value = 1

# Row 3 is brand new - clone A2's style onto A3 before writing the values.
$q1.Range("A2").Copy()
$q1.Range("A3").PasteSpecial(-4122)

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'006038"
$q1.Range("C3").Value = "大成景恒混合C"
$q1.Range("D3").Value = "'0.92"
$q1.Range("E3").Value = "'93.51"
$q1.Range("F3").Value = "'2.08"
$q1.Range("G3").Value = "'0.0191"
$q1.Range("H3").Value = 1

# ---------------------------------------------------------------------------
# 3) "总计" sheet: insert the new 2022-Q1 summary row, pushing the existing
#    2021-Q4 row down to row 3.
# ---------------------------------------------------------------------------

# Row 3 is new - clone A2's style onto A3 first.
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.01

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.07000000000000001
